# Weekly refresh of "Fruta / hortaliza" data: the per-row observation block
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Precio $/Kg -> columns D, J, K, L, M, P) is re-shuffled across the existing
# rows 2-36. Every other column (A, B, C, E-I, N, O, Q, R) stays untouched.
#
# Mapping: new row R gets the D/J/K/L/M/P values that old row $rowMap[R] had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 30;  3 = 12;  4 = 23;  5 = 19;  6 = 3;   7 = 10;  8 = 6;
    9 = 8;   10 = 33; 11 = 25; 12 = 32; 13 = 36; 14 = 28; 15 = 13;
    16 = 11; 17 = 15; 18 = 27; 19 = 29; 20 = 26; 21 = 31; 22 = 34;
    23 = 7;  24 = 14; 25 = 2;  26 = 22; 27 = 9;  28 = 4;  29 = 5;
    30 = 35; 31 = 18; 32 = 24; 33 = 17; 34 = 21; 35 = 16; 36 = 20
}

# Columns D, J, K, L, M, P as 1-based column indices (4, 10, 11, 12, 13, 16)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot the original values for every touched column/row before writing
# anything back, since source and destination rows overlap.
# NOTE: use Value2 (not Value) to read - Value getter is unreliable in this
# runtime and returns reflection metadata instead of the actual cell value.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 36; $r++) {
        $key = "$col,$r"
        $snapshot[$key] = $ws.Cells.Item($r, $col).Value2
    }
}

foreach ($col in $cols) {
    for ($r = 2; $r -le 36; $r++) {
        $srcRow = $rowMap[$r]
        $srcKey = "$col,$srcRow"
        $ws.Cells.Item($r, $col).Value2 = $snapshot[$srcKey]
    }
}
